$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 683.3333
$ws.Range("J12").Value = 550
$ws.Range("L12").Value = 550
$ws.Range("N12").Value = -890
$ws.Range("H38").Value = 1785.5714
$ws.Range("I38").Value = 1000
$ws.Range("J38").Value = 1916.5
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 5749.5
$ws.Range("M38").Value = -2628
$ws.Range("N38").Value = -6493.5
$ws.Range("H51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -10968
$ws.Range("H64").Value = 7250
$ws.Range("H67").Value = 7250
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H100").Value = 212.66667
$ws.Range("I100").Value = 212.66667
$ws.Range("K100").Value = 212.66667
$ws.Range("M100").Value = 328.33333
$ws.Range("H103").Value = 6066.25
$ws.Range("I103").Value = 10001
$ws.Range("K103").Value = 30003
$ws.Range("M103").Value = -29417
$ws.Range("H106").Value = 1955.6666
$ws.Range("I106").Value = 1955.6666
$ws.Range("K106").Value = 1955.6666
$ws.Range("M106").Value = -1324.6666
$ws.Range("H111").Value = 1791.4
$ws.Range("H112").Value = 2750.2144
$ws.Range("J112").Value = 2811.4443
$ws.Range("L112").Value = 8434.332900000001
$ws.Range("N112").Value = -10650.3329
$ws.Range("H113").Value = 2055.2
$ws.Range("I113").Value = 1856.5
$ws.Range("K113").Value = 1856.5
$ws.Range("M113").Value = 1397.5
$ws.Range("H132").Value = 8657.777
$ws.Range("I132").Value = 8899
$ws.Range("K132").Value = 26697
$ws.Range("M132").Value = -24167
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120
$ws.Range("H138").Value = 2863.2144
$ws.Range("I138").Value = 1098.6364
$ws.Range("J138").Value = 9333.333000000001
$ws.Range("K138").Value = 3295.9092
$ws.Range("L138").Value = 27999.999
$ws.Range("M138").Value = 1844.0908
$ws.Range("N138").Value = -38279.999
$ws.Range("H141").Value = 2338.3635
$ws.Range("I141").Value = 1873.2
$ws.Range("K141").Value = 5619.6
$ws.Range("M141").Value = -439.6000000000004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9091822
$ws.Range("I32").Value = 1003.8
$ws.Range("K32").Value = 1003.8
$ws.Range("M32").Value = -716.8
$ws.Range("H44").Value = 11865.315
$ws.Range("J44").Value = 11865.315
$ws.Range("L44").Value = 11865.315
$ws.Range("N44").Value = -12841.315
$ws.Range("H55").Value = 45666.332
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 60999.5
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 60999.5
$ws.Range("M55").Value = -14685
$ws.Range("N55").Value = -61629.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2757.5
$ws.Range("I20").Value = 3075.2
$ws.Range("K20").Value = 3075.2
$ws.Range("M20").Value = -2828.2
$ws.Range("H29").Value = 911.3333
$ws.Range("I29").Value = 948.6
$ws.Range("K29").Value = 948.6
$ws.Range("M29").Value = -659.6
$ws.Range("H36").Value = 838.4
$ws.Range("I36").Value = 838.4
$ws.Range("K36").Value = 838.4
$ws.Range("M36").Value = -304.4
$ws.Range("H81").Value = 53816.25
$ws.Range("J81").Value = 53816.25
$ws.Range("L81").Value = 53816.25
$ws.Range("N81").Value = -55938.25
$ws.Range("H84").Value = 53816.25
$ws.Range("J84").Value = 53816.25
$ws.Range("L84").Value = 161448.75
$ws.Range("N84").Value = -172056.75
$ws.Range("H134").Value = 8173
$ws.Range("I134").Value = 1744.625
$ws.Range("K134").Value = 5233.875
$ws.Range("M134").Value = -2698.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1067.25
$ws.Range("I105").Value = 756.3333
$ws.Range("K105").Value = 756.3333
$ws.Range("M105").Value = 990.6667
$ws.Range("H132").Value = 2470.6667
$ws.Range("I132").Value = 2081.7693
$ws.Range("K132").Value = 6245.3079
$ws.Range("M132").Value = -3715.3079
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.9
$ws.Range("I2").Value = 16.75
$ws.Range("K2").Value = 100.5
$ws.Range("M2").Value = 12.5
$ws.Range("H4").Value = 2400033.2
$ws.Range("I4").Value = 2400033.2
$ws.Range("K4").Value = 7200099.600000001
$ws.Range("M4").Value = -7199987.600000001
$ws.Range("H10").Value = 166.66667
$ws.Range("I10").Value = 166.66667
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500.00001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -361.00001
$ws.Range("N10").ClearContents()
$ws.Range("H16").Value = 100
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H51").Value = 2016
$ws.Range("I51").Value = 2021.5
$ws.Range("K51").Value = 6064.5
$ws.Range("M51").Value = -5604.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 198.46666
$ws.Range("I2").Value = 143.18182
$ws.Range("J2").Value = 350.5
$ws.Range("K2").Value = 143.18182
$ws.Range("L2").Value = 350.5
$ws.Range("M2").Value = -30.18181999999999
$ws.Range("N2").Value = -576.5
$ws.Range("H20").Value = 38094.855
$ws.Range("J20").Value = 38094.855
$ws.Range("L20").Value = 38094.855
$ws.Range("N20").Value = -38584.855
$ws.Range("H80").Value = 4302
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 4950.5
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 4950.5
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -6946.5
$ws.Range("H83").Value = 4302
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 4950.5
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 24752.5
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -34736.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("M3").Value = 12
$ws.Range("H12").Value = 344.77777
$ws.Range("J12").Value = 344.77777
$ws.Range("L12").Value = 344.77777
$ws.Range("N12").Value = -684.7777699999999
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 100
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = 70
$ws.Range("H46").Value = 6249
$ws.Range("I46").Value = 996.5
$ws.Range("K46").Value = 996.5
$ws.Range("M46").Value = -808.5
$ws.Range("H76").Value = 26525.334
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 26525.334
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 26525.334
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -27201.334
$ws.Range("H79").Value = 26525.334
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 26525.334
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 26525.334
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -28865.334
$ws.Range("H132").Value = 7332.3335
$ws.Range("I132").Value = 8399.799999999999
$ws.Range("K132").Value = 25199.4
$ws.Range("M132").Value = -22669.4
$ws.Range("H135").Value = 78666.336
$ws.Range("J135").Value = 78666.336
$ws.Range("L135").Value = 78666.336
$ws.Range("N135").Value = -88806.336
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1006.6667
$ws.Range("I2").Value = 160
$ws.Range("J2").Value = 2700
$ws.Range("K2").Value = 160
$ws.Range("L2").Value = 2700
$ws.Range("M2").Value = -48
$ws.Range("N2").Value = -2924
$ws.Range("H4").Value = 23900
$ws.Range("I4").Value = 35833.332
$ws.Range("K4").Value = 35833.332
$ws.Range("M4").Value = -35720.332
$ws.Range("H5").Value = 12300000
$ws.Range("I5").Value = 13020000
$ws.Range("K5").Value = 13020000
$ws.Range("M5").Value = -13019888
$ws.Range("H7").Value = 1600
$ws.Range("I7").Value = 1400
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1287
$ws.Range("N7").Value = -2226
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H136").Value = 3142.037
$ws.Range("I136").Value = 2155.375
$ws.Range("K136").Value = 6466.125
$ws.Range("M136").Value = -3916.125
